# Fruta / hortaliza, semanal
# Insert a new weekly record above the existing row 328 on the single
# worksheet, shifting all subsequent rows (328-349) down by one (to
# 329-350) and filling the newly opened row 328 with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 328; Excel shifts rows 328..349
# down to 329..350, preserving their formatting/styles (including the
# date-format style on column D).
$ws.Rows("328:328").Insert()

# Populate the newly inserted row 328 with the new weekly data point.
$ws.Range("A328").Value = 4
$ws.Range("B328").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C328").Value = "Los Lagos"
$ws.Range("D328").Value = 45013
$ws.Range("E328").Value = 10
$ws.Range("F328").Value = "Fruta"
$ws.Range("G328").Value = 100108
$ws.Range("H328").Value = "Tropicales y subtropicales"
$ws.Range("I328").Value = 100108002
$ws.Range("J328").Value = "Mango"
$ws.Range("K328").Value = "Sin especificar"
$ws.Range("L328").Value = "Primera"
$ws.Range("M328").Value = 200
$ws.Range("N328").Value = 8000
$ws.Range("O328").Value = 8500
$ws.Range("P328").Value = 8250
$ws.Range("Q328").Value = "`$/bandeja 4 kilos"
$ws.Range("R328").Value = "Perú"
$ws.Range("S328").Value = 2062
$ws.Range("T328").Value = 4
